$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right before row 221, shifting the existing
# rows 221-288 down to 223-290 (dimension grows from A1:T288 to A1:T290).
$ws.Range("A221:T222").EntireRow.Insert()

# New row 221: Navel Late / Primera (Region de O'Higgins)
$ws.Range("A221").Value = 5
$ws.Range("B221").Value = "Macroferia Regional de Talca"
$ws.Range("C221").Value = "Maule"
$ws.Range("D221").Value = 44468
$ws.Range("E221").Value = 7
$ws.Range("F221").Value = "Fruta"
$ws.Range("G221").Value = 100102
$ws.Range("H221").Value = "Cítricos"
$ws.Range("I221").Value = 100102005
$ws.Range("J221").Value = "Naranja"
$ws.Range("K221").Value = "Navel Late"
$ws.Range("L221").Value = "Primera"
$ws.Range("M221").Value = 300
$ws.Range("N221").Value = 7000
$ws.Range("O221").Value = 7000
$ws.Range("P221").Value = 7000
$ws.Range("Q221").Value = "$/bandeja 15 kilos granel"
$ws.Range("R221").Value = "Región de O'Higgins"
$ws.Range("S221").Value = 467
$ws.Range("T221").Value = 15

# New row 222: Washington parent / Primera (Provincia de Quillota)
$ws.Range("A222").Value = 5
$ws.Range("B222").Value = "Macroferia Regional de Talca"
$ws.Range("C222").Value = "Maule"
$ws.Range("D222").Value = 44468
$ws.Range("E222").Value = 7
$ws.Range("F222").Value = "Fruta"
$ws.Range("G222").Value = 100102
$ws.Range("H222").Value = "Cítricos"
$ws.Range("I222").Value = 100102005
$ws.Range("J222").Value = "Naranja"
$ws.Range("K222").Value = "Washington parent"
$ws.Range("L222").Value = "Primera"
$ws.Range("M222").Value = 230
$ws.Range("N222").Value = 5000
$ws.Range("O222").Value = 5000
$ws.Range("P222").Value = 5000
$ws.Range("Q222").Value = "$/bandeja 15 kilos granel"
$ws.Range("R222").Value = "Provincia de Quillota"
$ws.Range("S222").Value = 333
$ws.Range("T222").Value = 15

# Match the date-formatted style used by the rest of column D.
$ws.Range("D221").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D222").NumberFormat = "YYYY-MM-DD HH:MM:SS"
